# "Generate Report for Archive"
# The localization status report is regenerated: items that were previously
# "Ready for handoff" have moved on to "In Translation", and (since that
# label is shorter) the Status-ish columns that display it are re-sized to
# fit the new content.

$wb = $excel.ActiveWorkbook

# Update the status text everywhere it appears (Overview!E:F and the
# per-locale sheets' Status column) in one shot, across every worksheet.
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# Re-fit the columns that held the old, longer "Ready for handoff" text so
# they match the narrower width needed for "In Translation".
$target = 13.4101845877511
$newColumnWidth = $target - 0.8333333333333334

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth  # column E: zh-cn
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth  # column F: de-de

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth      # column C: Status

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth      # column C: Status
